# Apply updated crypto price/volume figures per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.258.05"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.503.40"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'584.52"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'134.89"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").Value = "3.504.70"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'7.14"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").Value = "4.098.72"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "3.504.90"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'26.41"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").Value = "64.266.74"
$ws.Range("D19").Value = "'9.77"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'13.87"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").Value = "'383.91"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "3.640.27"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'73.95"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "'7.54"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "3.522.56"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'5.33"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("D41").Value = "'164.45"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("D42").Value = "'0.0783"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").Value = "'0.810"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "'26.07"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'41.85"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'4.42"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "2.474.92"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'0.922"
$ws.Range("E51").Value = "  +1.11%  "
